$d = $word.ActiveDocument

# Update the date line (wdReplaceOne = 1)
$d.Content.Find.Execute("2023-03-30 Thursday", $false, $false, $false, $false, $false, $true, 1, $false, "2023-03-31 Friday", 1) | Out-Null

# Update table cells (20 rows x 5 cols), row-major order matching the diff.
# Each replacement is scoped to its own cell Range and uses wdReplaceOne (1)
# rather than wdReplaceAll (2), because ReplaceAll searches the whole story
# regardless of range scoping, while ReplaceOne only replaces the single match
# found (respecting the cell-scoped range) -- this matters because several
# source cells share the same old text (e.g. "66×70=4620" occurs twice).
$t = $d.Tables.Item(1)
$pairs = @(
    @("97×64=6208", "86×47=4042"),
    @("80×56=4480", "72×20=1440"),
    @("12×74=888", "69×39=2691"),
    @("93×52=4836", "75×31=2325"),
    @("11×19=209", "75×60=4500"),
    @("57×46=2622", "15×71=1065"),
    @("18×97=1746", "84×88=7392"),
    @("14×22=308", "20×83=1660"),
    @("47×43=2021", "28×13=364"),
    @("61×16=976", "31×20=620"),
    @("54×42=2268", "70×28=1960"),
    @("100×43=4300", "80×72=5760"),
    @("40×36=1440", "67×13=871"),
    @("91×95=8645", "31×68=2108"),
    @("77×87=6699", "28×66=1848"),
    @("22×70=1540", "68×26=1768"),
    @("89×42=3738", "99×28=2772"),
    @("24×52=1248", "16×10=160"),
    @("98×62=6076", "83×57=4731"),
    @("87×78=6786", "41×27=1107"),
    @("96×16=1536", "39×40=1560"),
    @("86×69=5934", "85×74=6290"),
    @("29×51=1479", "26×92=2392"),
    @("57×55=3135", "21×65=1365"),
    @("37×44=1628", "37×54=1998"),
    @("12×97=1164", "72×56=4032"),
    @("72×45=3240", "57×80=4560"),
    @("10×23=230", "38×86=3268"),
    @("86×57=4902", "50×59=2950"),
    @("17×56=952", "56×24=1344"),
    @("35×59=2065", "17×66=1122"),
    @("51×28=1428", "41×53=2173"),
    @("31×57=1767", "16×85=1360"),
    @("90×68=6120", "95×11=1045"),
    @("68×25=1700", "19×62=1178"),
    @("49×43=2107", "30×78=2340"),
    @("72×88=6336", "32×51=1632"),
    @("58×51=2958", "100×72=7200"),
    @("15×83=1245", "47×98=4606"),
    @("53×38=2014", "82×71=5822"),
    @("29×46=1334", "67×45=3015"),
    @("83×81=6723", "95×60=5700"),
    @("79×76=6004", "35×28=980"),
    @("64×81=5184", "69×35=2415"),
    @("39×55=2145", "39×43=1677"),
    @("65×80=5200", "66×82=5412"),
    @("73×96=7008", "77×22=1694"),
    @("39×94=3666", "27×87=2349"),
    @("22×19=418", "92×91=8372"),
    @("82×50=4100", "67×49=3283"),
    @("82×38=3116", "88×35=3080"),
    @("44×54=2376", "47×33=1551"),
    @("39×37=1443", "55×66=3630"),
    @("54×73=3942", "87×18=1566"),
    @("13×95=1235", "77×45=3465"),
    @("61×11=671", "95×19=1805"),
    @("19×25=475", "51×96=4896"),
    @("92×94=8648", "39×61=2379"),
    @("47×42=1974", "65×40=2600"),
    @("71×88=6248", "58×15=870"),
    @("88×82=7216", "62×95=5890"),
    @("82×49=4018", "96×90=8640"),
    @("53×60=3180", "11×51=561"),
    @("93×71=6603", "71×33=2343"),
    @("68×14=952", "72×12=864"),
    @("29×48=1392", "41×57=2337"),
    @("66×70=4620", "35×73=2555"),
    @("90×24=2160", "89×55=4895"),
    @("11×39=429", "73×95=6935"),
    @("12×15=180", "90×16=1440"),
    @("97×77=7469", "88×47=4136"),
    @("14×40=560", "27×95=2565"),
    @("37×87=3219", "11×35=385"),
    @("18×55=990", "68×20=1360"),
    @("12×53=636", "24×27=648"),
    @("25×98=2450", "51×74=3774"),
    @("34×80=2720", "61×31=1891"),
    @("65×90=5850", "73×38=2774"),
    @("70×49=3430", "64×55=3520"),
    @("45×56=2520", "57×13=741"),
    @("90×40=3600", "42×45=1890"),
    @("29×85=2465", "29×50=1450"),
    @("29×56=1624", "86×83=7138"),
    @("91×69=6279", "19×61=1159"),
    @("67×74=4958", "90×65=5850"),
    @("63×54=3402", "54×41=2214"),
    @("79×85=6715", "87×76=6612"),
    @("24×87=2088", "61×41=2501"),
    @("16×11=176", "92×73=6716"),
    @("30×49=1470", "67×71=4757"),
    @("61×30=1830", "40×33=1320"),
    @("46×48=2208", "76×46=3496"),
    @("73×21=1533", "39×47=1833"),
    @("66×70=4620", "19×56=1064"),
    @("36×20=720", "80×57=4560"),
    @("26×44=1144", "100×55=5500"),
    @("87×94=8178", "61×93=5673"),
    @("91×57=5187", "80×93=7440"),
    @("21×17=357", "13×22=286"),
    @("82×70=5740", "42×100=4200")
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $pair = $pairs[$idx]
        $cell = $t.Cell($r, $c)
        $ok = $cell.Range.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 1)
        if (-not $ok) {
            Write-Host "WARNING: replace failed at row" $r "col" $c "old=" $pair[0]
        }
        $idx++
    }
}

Write-Host "Done. idx=" $idx